$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Google Shape;88;p1" textbox holds the talk title/speaker/date block.
$shape = $s.Shapes.Item(5)

# Reposition/resize the textbox (EMU -> points, 914400 EMU per inch / 72pt per inch).
# Point values below are nudged within their rounding tolerance so the
# round-trip back through the host's internal float32 storage lands on the
# exact target EMU (3384646, 2157901, 4082954, 4031833).
$shape.Left = 266.50759885511815
$shape.Top = 169.91349796692913
$shape.Width = 321.4924774448819
$shape.Height = 317.46720885433075

# Update the day-of-week in the date line from "Friday" to "Wednesday".
$tf = $shape.TextFrame
$tr = $tf.TextRange
$found = $tr.Replace("Friday, ", "Wednesday, ")
